# Trading update: 2026-02-17 08:33:02
# Appends the latest MarketMaking trade (trade #38, still OPEN) as a new
# row at the bottom of both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $row = 39

    $ws.Cells.Item($row, 1).Value = 38

    # Dates like "2026-02-17" are auto-recognized as date serials by
    # Excel, but the sheet stores them as plain text - force a text
    # format before writing, then drop back to the Normal style so no
    # leftover number formatting is left on the cell.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "08:33:00"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.67
    # Column G (Exit Price) is left blank - the trade is still OPEN.
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 99.58598934440597
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    # Column P (Exit Reason) is left blank - the trade is still OPEN.
    $ws.Cells.Item($row, 17).Value = 0
}
